$wb = $excel.ActiveWorkbook

# Add a new worksheet named "Sayfa1" after the last existing sheet
$lastIndex = $wb.Worksheets.Count
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($lastIndex))
$newSheet.Name = "Sayfa1"

# Fill in the data
$newSheet.Range("A1").Value = "anam"
$newSheet.Range("B1").Value = 1
$newSheet.Range("A2").Value = "babam"
$newSheet.Range("B2").Value = 2
$newSheet.Range("A3").Value = "cnm"
$newSheet.Range("B3").Value = 2
$newSheet.Range("A4").Value = "benim"
$newSheet.Range("B4").Value = 3

# Select the new sheet and set the active cell like the original edit
$newSheet.Activate()
$newSheet.Range("D7").Select()
